# Update master to output generated at c986bee
$d = $word.ActiveDocument

# 1. Update the date line at the top of the document.
$d.Content.Find.Execute("2025-01-01 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-01-02 Thursday", 2) | Out-Null

# 2. Update the answer table. Rows 2/3/4 (and the 3 blank rows following each
#    data row) are untouched; only the five data rows change their values.
$t = $d.Tables.Item(1)

$rowValues = @{
    1  = @("22÷7=3, 1", "53÷7=7, 4", "53÷4=13, 1", "33÷7=4, 5", "92÷8=11, 4")
    5  = @("98÷4=24, 2", "59÷3=19, 2", "78÷6=13, 0", "13÷4=3, 1", "31÷4=7, 3")
    9  = @("38÷9=4, 2", "71÷3=23, 2", "44÷5=8, 4", "43÷5=8, 3", "85÷3=28, 1")
    13 = @("42÷4=10, 2", "61÷2=30, 1", "77÷4=19, 1", "91÷3=30, 1", "73÷2=36, 1")
    17 = @("53÷4=13, 1", "17÷6=2, 5", "83÷7=11, 6", "78÷5=15, 3", "83÷7=11, 6")
}

foreach ($rowIndex in $rowValues.Keys) {
    $values = $rowValues[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]
    }
}
